$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 73; existing rows 73-125 shift down to 74-126.
$ws.Rows.Item(73).Insert()

$newRow = 73
$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = (Get-Date -Year 2022 -Month 7 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112022
$ws.Cells.Item($newRow, 7).Value = "Arveja Verde"
$ws.Cells.Item($newRow, 8).Value = "Perfection"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 50
$ws.Cells.Item($newRow, 11).Value = 45000
$ws.Cells.Item($newRow, 12).Value = 45000
$ws.Cells.Item($newRow, 13).Value = 45000
$ws.Cells.Item($newRow, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item($newRow, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($newRow, 16).Value = 1800
$ws.Cells.Item($newRow, 17).Value = 25
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
